# Rebuild the "Export" sheet data (Conta / Nome / Saldo) to match the new
# source export: several accounts were added, removed and reordered, and
# several balances were refreshed. We rewrite the whole used range in one
# shot via a 2D array for reliability, rather than patching individual rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 215,3
$data[0,0] = "Conta"
$data[0,1] = "Nome"
$data[0,2] = "Saldo"
$data[1,0] = "005646524"
$data[1,1] = "EVANGELINA"
$data[1,2] = 400000
$data[2,0] = "004384167"
$data[2,1] = "DOUGLAS"
$data[2,2] = 3668.89
$data[3,0] = "000806386"
$data[3,1] = "FERNANDA"
$data[3,2] = 3153.98
$data[4,0] = "004329030"
$data[4,1] = "DANIELA"
$data[4,2] = 3088.27
$data[5,0] = "004207278"
$data[5,1] = "CESAR"
$data[5,2] = 2000
$data[6,0] = "004368468"
$data[6,1] = "AHMAD"
$data[6,2] = 1827.75
$data[7,0] = "005591536"
$data[7,1] = "GUSTAVO"
$data[7,2] = 1283.6
$data[8,0] = "002697806"
$data[8,1] = "CLAUDIA"
$data[8,2] = 1156.52
$data[9,0] = "001651617"
$data[9,1] = "MIRELLA"
$data[9,2] = 1010.5
$data[10,0] = "008119302"
$data[10,1] = "VITOR"
$data[10,2] = 1003.77
$data[11,0] = "008110684"
$data[11,1] = "EDVAL"
$data[11,2] = 935.16
$data[12,0] = "005883672"
$data[12,1] = "FLK"
$data[12,2] = 908.42
$data[13,0] = "004392159"
$data[13,1] = "RODRIGO"
$data[13,2] = 900.21
$data[14,0] = "004471893"
$data[14,1] = "PAULA"
$data[14,2] = 887.53
$data[15,0] = "004212581"
$data[15,1] = "MARIA"
$data[15,2] = 880.8
$data[16,0] = "001761119"
$data[16,1] = "BLUEMETRIX"
$data[16,2] = 868.61
$data[17,0] = "004458624"
$data[17,1] = "PEDRO"
$data[17,2] = 849.33
$data[18,0] = "004381359"
$data[18,1] = "PEDRO"
$data[18,2] = 837.95
$data[19,0] = "005685353"
$data[19,1] = "CARLOS"
$data[19,2] = 827.8
$data[20,0] = "004322719"
$data[20,1] = "GISELA"
$data[20,2] = 753.02
$data[21,0] = "004693308"
$data[21,1] = "LAURA"
$data[21,2] = 706
$data[22,0] = "005440756"
$data[22,1] = "VALERIA"
$data[22,2] = 705.15
$data[23,0] = "008026930"
$data[23,1] = "JOAO"
$data[23,2] = 700
$data[24,0] = "004237325"
$data[24,1] = "RICARDO"
$data[24,2] = 617.14
$data[25,0] = "004514241"
$data[25,1] = "ANDRE"
$data[25,2] = 586.4
$data[26,0] = "004975924"
$data[26,1] = "SERGIO"
$data[26,2] = 574.44
$data[27,0] = "004352384"
$data[27,1] = "BRASFORT"
$data[27,2] = 564.58
$data[28,0] = "005152037"
$data[28,1] = "RODRIGO"
$data[28,2] = 562.4
$data[29,0] = "003512801"
$data[29,1] = "LAIS"
$data[29,2] = 538.94
$data[30,0] = "004377713"
$data[30,1] = "DANIELI"
$data[30,2] = 496.98
$data[31,0] = "008115273"
$data[31,1] = "NILSON"
$data[31,2] = 492.2
$data[32,0] = "005044389"
$data[32,1] = "CLAUDIA"
$data[32,2] = 461.89
$data[33,0] = "005637820"
$data[33,1] = "GUILHERME"
$data[33,2] = 459.95
$data[34,0] = "004207374"
$data[34,1] = "ANGELICA"
$data[34,2] = 450.48
$data[35,0] = "005142592"
$data[35,1] = "ALBERTO"
$data[35,2] = 450
$data[36,0] = "004432579"
$data[36,1] = "ANA"
$data[36,2] = 446.18
$data[37,0] = "005924958"
$data[37,1] = "TIAGO"
$data[37,2] = 438.4
$data[38,0] = "004436055"
$data[38,1] = "MARCO"
$data[38,2] = 365.23
$data[39,0] = "003553997"
$data[39,1] = "MIRELLA"
$data[39,2] = 341.62
$data[40,0] = "004556853"
$data[40,1] = "MARCEL"
$data[40,2] = 336.04
$data[41,0] = "004424761"
$data[41,1] = "PEDRO"
$data[41,2] = 330
$data[42,0] = "004413523"
$data[42,1] = "ROSANE"
$data[42,2] = 304.16
$data[43,0] = "008054285"
$data[43,1] = "IGOR"
$data[43,2] = 277.63
$data[44,0] = "004424671"
$data[44,1] = "LUISA"
$data[44,2] = 250
$data[45,0] = "004424714"
$data[45,1] = "HELENA"
$data[45,2] = 250
$data[46,0] = "003249855"
$data[46,1] = "MARINA"
$data[46,2] = 237.12
$data[47,0] = "005295509"
$data[47,1] = "BHRUNA"
$data[47,2] = 223.02
$data[48,0] = "004467884"
$data[48,1] = "ANA"
$data[48,2] = 193.66
$data[49,0] = "004487016"
$data[49,1] = "ROGERIO"
$data[49,2] = 184.88
$data[50,0] = "004207184"
$data[50,1] = "CRISTINA"
$data[50,2] = 177.33
$data[51,0] = "004208447"
$data[51,1] = "LEILA"
$data[51,2] = 161.7
$data[52,0] = "004508526"
$data[52,1] = "CASSIO"
$data[52,2] = 152.51
$data[53,0] = "005142611"
$data[53,1] = "GUILHERME"
$data[53,2] = 134.69
$data[54,0] = "001719494"
$data[54,1] = "LUIS"
$data[54,2] = 106.95
$data[55,0] = "005141215"
$data[55,1] = "KARINA"
$data[55,2] = 100
$data[56,0] = "008071998"
$data[56,1] = "ISADORA"
$data[56,2] = 100
$data[57,0] = "004211911"
$data[57,1] = "ZENILDA"
$data[57,2] = 96
$data[58,0] = "004239387"
$data[58,1] = "LUIZ"
$data[58,2] = 95.08
$data[59,0] = "004431591"
$data[59,1] = "MARIO"
$data[59,2] = 94.24
$data[60,0] = "004218542"
$data[60,1] = "JOSE"
$data[60,2] = 93.52
$data[61,0] = "008035153"
$data[61,1] = "CLAUDIO"
$data[61,2] = 92.33
$data[62,0] = "004749928"
$data[62,1] = "NADY"
$data[62,2] = 92.24
$data[63,0] = "003115072"
$data[63,1] = "VICTOR"
$data[63,2] = 91.24
$data[64,0] = "004340984"
$data[64,1] = "RENATA"
$data[64,2] = 90.99
$data[65,0] = "005268516"
$data[65,1] = "LUIS"
$data[65,2] = 89.34
$data[66,0] = "004425965"
$data[66,1] = "CAROLLINA"
$data[66,2] = 87.87
$data[67,0] = "004466350"
$data[67,1] = "RAQUEL"
$data[67,2] = 87.36
$data[68,0] = "005266369"
$data[68,1] = "EG"
$data[68,2] = 87.08
$data[69,0] = "004383268"
$data[69,1] = "LAURA"
$data[69,2] = 86.58
$data[70,0] = "004384258"
$data[70,1] = "PAULA"
$data[70,2] = 86.58
$data[71,0] = "004212132"
$data[71,1] = "JOAO"
$data[71,2] = 86.38
$data[72,0] = "004809902"
$data[72,1] = "PEDRO"
$data[72,2] = 85.9
$data[73,0] = "004536602"
$data[73,1] = "TATIANY"
$data[73,2] = 85.74
$data[74,0] = "005312963"
$data[74,1] = "ALAN"
$data[74,2] = 84.58
$data[75,0] = "004260002"
$data[75,1] = "ERICA"
$data[75,2] = 84.52
$data[76,0] = "004752615"
$data[76,1] = "LUZIMAR"
$data[76,2] = 84.49
$data[77,0] = "008149996"
$data[77,1] = "CAMILA"
$data[77,2] = 84.41
$data[78,0] = "004332103"
$data[78,1] = "JOSE"
$data[78,2] = 80.96
$data[79,0] = "004261201"
$data[79,1] = "ANA"
$data[79,2] = 79.12
$data[80,0] = "004272426"
$data[80,1] = "RODRIGO"
$data[80,2] = 78.52
$data[81,0] = "005186167"
$data[81,1] = "ANDREA"
$data[81,2] = 77.77
$data[82,0] = "003836362"
$data[82,1] = "ISABELLA"
$data[82,2] = 77.34
$data[83,0] = "005701765"
$data[83,1] = "F"
$data[83,2] = 75.98
$data[84,0] = "005206566"
$data[84,1] = "LEVI"
$data[84,2] = 71.36
$data[85,0] = "004563252"
$data[85,1] = "FERNANDO"
$data[85,2] = 70.58
$data[86,0] = "005068961"
$data[86,1] = "JORGE"
$data[86,2] = 70.16
$data[87,0] = "004290978"
$data[87,1] = "LARISSA"
$data[87,2] = 69
$data[88,0] = "004472760"
$data[88,1] = "SANDRA"
$data[88,2] = 68.77
$data[89,0] = "004691225"
$data[89,1] = "ANNA"
$data[89,2] = 67.55
$data[90,0] = "005018038"
$data[90,1] = "ELAINE"
$data[90,2] = 67.28
$data[91,0] = "004508159"
$data[91,1] = "FELIPE"
$data[91,2] = 66.87
$data[92,0] = "005558076"
$data[92,1] = "ALEXANDRE"
$data[92,2] = 65.01
$data[93,0] = "004752519"
$data[93,1] = "MARCUS"
$data[93,2] = 63.26
$data[94,0] = "004974089"
$data[94,1] = "CELIA"
$data[94,2] = 59.36
$data[95,0] = "004877741"
$data[95,1] = "LUIZ"
$data[95,2] = 59.12
$data[96,0] = "005685089"
$data[96,1] = "CARNEIRO"
$data[96,2] = 58.58
$data[97,0] = "005880251"
$data[97,1] = "LUIZ"
$data[97,2] = 52.3
$data[98,0] = "004400640"
$data[98,1] = "FELIPE"
$data[98,2] = 51.44
$data[99,0] = "008054713"
$data[99,1] = "MODULAR"
$data[99,2] = 51.43
$data[100,0] = "004994036"
$data[100,1] = "BALTASAR"
$data[100,2] = 50.65
$data[101,0] = "005076418"
$data[101,1] = "LEONARDO"
$data[101,2] = 50.23
$data[102,0] = "004477812"
$data[102,1] = "DIEGO"
$data[102,2] = 48.45
$data[103,0] = "003497496"
$data[103,1] = "ELISANDRA"
$data[103,2] = 48.02
$data[104,0] = "004546050"
$data[104,1] = "LUIS"
$data[104,2] = 47.04
$data[105,0] = "004971448"
$data[105,1] = "CLOVIS"
$data[105,2] = 45.87
$data[106,0] = "008115927"
$data[106,1] = "ARI"
$data[106,2] = 44.96
$data[107,0] = "001731007"
$data[107,1] = "GUILHERME"
$data[107,2] = 44.73
$data[108,0] = "005103059"
$data[108,1] = "WALQUIRIA"
$data[108,2] = 41.88
$data[109,0] = "004452507"
$data[109,1] = "DANIELA"
$data[109,2] = 41.43
$data[110,0] = "005514036"
$data[110,1] = "ANA"
$data[110,2] = 41.34
$data[111,0] = "008069255"
$data[111,1] = "ANGELA"
$data[111,2] = 40.91
$data[112,0] = "004242237"
$data[112,1] = "MARIAH"
$data[112,2] = 39.99
$data[113,0] = "004381194"
$data[113,1] = "ALINNE"
$data[113,2] = 39.91
$data[114,0] = "005245032"
$data[114,1] = "ROSA"
$data[114,2] = 39.91
$data[115,0] = "004238164"
$data[115,1] = "DANIELA"
$data[115,2] = 39.09
$data[116,0] = "004920447"
$data[116,1] = "MARILIA"
$data[116,2] = 38.63
$data[117,0] = "004481463"
$data[117,1] = "MARA"
$data[117,2] = 37.31
$data[118,0] = "004584517"
$data[118,1] = "CAIO"
$data[118,2] = 36.51
$data[119,0] = "004806286"
$data[119,1] = "VERA"
$data[119,2] = 35.77
$data[120,0] = "004452912"
$data[120,1] = "BRUNO"
$data[120,2] = 35.75
$data[121,0] = "005079458"
$data[121,1] = "JONAS"
$data[121,2] = 35.4
$data[122,0] = "004981655"
$data[122,1] = "WILLIAM"
$data[122,2] = 34.5
$data[123,0] = "004398174"
$data[123,1] = "DANIELE"
$data[123,2] = 34.44
$data[124,0] = "004332207"
$data[124,1] = "IRACY"
$data[124,2] = 34.03
$data[125,0] = "004213139"
$data[125,1] = "LEONARDO"
$data[125,2] = 33.39
$data[126,0] = "004230529"
$data[126,1] = "LAIS"
$data[126,2] = 31.08
$data[127,0] = "005927101"
$data[127,1] = "SIMONE"
$data[127,2] = 30
$data[128,0] = "005305965"
$data[128,1] = "SIDMAR"
$data[128,2] = 28.77
$data[129,0] = "004377415"
$data[129,1] = "ANGELA"
$data[129,2] = 28.73
$data[130,0] = "004404724"
$data[130,1] = "LEANDRO"
$data[130,2] = 28.45
$data[131,0] = "004228456"
$data[131,1] = "FLASH"
$data[131,2] = 27.46
$data[132,0] = "004472076"
$data[132,1] = "RUBENS"
$data[132,2] = 27.42
$data[133,0] = "004350197"
$data[133,1] = "GISELA"
$data[133,2] = 25.08
$data[134,0] = "008002502"
$data[134,1] = "JORGEANA"
$data[134,2] = 24.87
$data[135,0] = "004405476"
$data[135,1] = "MARIANA"
$data[135,2] = 24.52
$data[136,0] = "002064834"
$data[136,1] = "RAFAELA"
$data[136,2] = 23.64
$data[137,0] = "005255637"
$data[137,1] = "PATRICIA"
$data[137,2] = 23.19
$data[138,0] = "000827730"
$data[138,1] = "LUCIANA"
$data[138,2] = 21.82
$data[139,0] = "004371857"
$data[139,1] = "NAZARETH"
$data[139,2] = 21.52
$data[140,0] = "004388077"
$data[140,1] = "WLADMIR"
$data[140,2] = 20.89
$data[141,0] = "004360431"
$data[141,1] = "CARLOS"
$data[141,2] = 20.81
$data[142,0] = "004214604"
$data[142,1] = "MARIA"
$data[142,2] = 20.75
$data[143,0] = "004381415"
$data[143,1] = "JOAO"
$data[143,2] = 20.22
$data[144,0] = "004328934"
$data[144,1] = "VALERIA"
$data[144,2] = 19.61
$data[145,0] = "008026942"
$data[145,1] = "THOMAS"
$data[145,2] = 19.38
$data[146,0] = "005324981"
$data[146,1] = "JO"
$data[146,2] = 19.01
$data[147,0] = "004204255"
$data[147,1] = "AMADO"
$data[147,2] = 18.77
$data[148,0] = "004907688"
$data[148,1] = "HEITOR"
$data[148,2] = 18.23
$data[149,0] = "005981575"
$data[149,1] = "GLAUCIANE"
$data[149,2] = 16.7
$data[150,0] = "002894447"
$data[150,1] = "JOAO"
$data[150,2] = 16.48
$data[151,0] = "004264780"
$data[151,1] = "MARCELO"
$data[151,2] = 14.9
$data[152,0] = "004422594"
$data[152,1] = "WANDIR"
$data[152,2] = 14.67
$data[153,0] = "008032257"
$data[153,1] = "SARA"
$data[153,2] = 14.53
$data[154,0] = "004505474"
$data[154,1] = "RICARDO"
$data[154,2] = 13.23
$data[155,0] = "004498637"
$data[155,1] = "TIAGO"
$data[155,2] = 11.08
$data[156,0] = "005374916"
$data[156,1] = "MARCO"
$data[156,2] = 10.89
$data[157,0] = "005173958"
$data[157,1] = "VENIA"
$data[157,2] = 10.62
$data[158,0] = "004216298"
$data[158,1] = "FLORDELIZ"
$data[158,2] = 9.8
$data[159,0] = "005077648"
$data[159,1] = "DUNAS"
$data[159,2] = 8.68
$data[160,0] = "004289402"
$data[160,1] = "LARISSA"
$data[160,2] = 8.1
$data[161,0] = "008012870"
$data[161,1] = "ANA"
$data[161,2] = 7.7
$data[162,0] = "008004995"
$data[162,1] = "JOSE"
$data[162,2] = 7.14
$data[163,0] = "004530494"
$data[163,1] = "ROSANGELA"
$data[163,2] = 7
$data[164,0] = "008013889"
$data[164,1] = "CAROLINA"
$data[164,2] = 6
$data[165,0] = "004224405"
$data[165,1] = "MILA"
$data[165,2] = 5.88
$data[166,0] = "004448501"
$data[166,1] = "JOAO"
$data[166,2] = 5.55
$data[167,0] = "008032413"
$data[167,1] = "VICTOR"
$data[167,2] = 5.53
$data[168,0] = "005198093"
$data[168,1] = "ANA"
$data[168,2] = 4.66
$data[169,0] = "000834301"
$data[169,1] = "MARCUS"
$data[169,2] = 4.4
$data[170,0] = "004221638"
$data[170,1] = "CAROLINE"
$data[170,2] = 3.54
$data[171,0] = "004382902"
$data[171,1] = "LEILA"
$data[171,2] = 3.24
$data[172,0] = "004488571"
$data[172,1] = "CARLOS"
$data[172,2] = 1.74
$data[173,0] = "005022526"
$data[173,1] = "ALEXANDRE"
$data[173,2] = 1.7
$data[174,0] = "004754056"
$data[174,1] = "BRUNO"
$data[174,2] = 1.56
$data[175,0] = "004359408"
$data[175,1] = "HEPTA"
$data[175,2] = 1.55
$data[176,0] = "000431814"
$data[176,1] = "GUILHERME"
$data[176,2] = 1.1
$data[177,0] = "004360430"
$data[177,1] = "VIOMAR"
$data[177,2] = 1
$data[178,0] = "004486497"
$data[178,1] = "ELENA"
$data[178,2] = 0.96
$data[179,0] = "004115403"
$data[179,1] = "HEBERT"
$data[179,2] = 0.88
$data[180,0] = "005660155"
$data[180,1] = "CAROLINA"
$data[180,2] = 0.85
$data[181,0] = "004223502"
$data[181,1] = "BRUNA"
$data[181,2] = 0.78
$data[182,0] = "004862746"
$data[182,1] = "CESAR"
$data[182,2] = 0.71
$data[183,0] = "004587511"
$data[183,1] = "CARLOS"
$data[183,2] = 0.69
$data[184,0] = "004473942"
$data[184,1] = "DAIANNE"
$data[184,2] = 0.62
$data[185,0] = "003894173"
$data[185,1] = "ANDREA"
$data[185,2] = 0.48
$data[186,0] = "005645211"
$data[186,1] = "AGUINALDO"
$data[186,2] = 0.45
$data[187,0] = "004565108"
$data[187,1] = "GUSTAVO"
$data[187,2] = 0.42
$data[188,0] = "004453302"
$data[188,1] = "ISABELLA"
$data[188,2] = 0.39
$data[189,0] = "005749972"
$data[189,1] = "ALESSANDRA"
$data[189,2] = 0.37
$data[190,0] = "004278033"
$data[190,1] = "DAISY"
$data[190,2] = 0.21
$data[191,0] = "001759765"
$data[191,1] = "NATAL"
$data[191,2] = 0.19
$data[192,0] = "004432455"
$data[192,1] = "LUCIANA"
$data[192,2] = 0.17
$data[193,0] = "002694089"
$data[193,1] = "VITOR"
$data[193,2] = 0.16
$data[194,0] = "004357159"
$data[194,1] = "JOAO"
$data[194,2] = 0.15
$data[195,0] = "004320840"
$data[195,1] = "NATALIA"
$data[195,2] = 0.14
$data[196,0] = "001000288"
$data[196,1] = "ISABELLA"
$data[196,2] = 0.13
$data[197,0] = "005530256"
$data[197,1] = "CAROLINA"
$data[197,2] = 0.11
$data[198,0] = "004451996"
$data[198,1] = "ADRIANO"
$data[198,2] = 0.09
$data[199,0] = "005047946"
$data[199,1] = "GABRIEL"
$data[199,2] = 0.09
$data[200,0] = "005075382"
$data[200,1] = "NAYARA"
$data[200,2] = 0.09
$data[201,0] = "004335251"
$data[201,1] = "EDMUNDO"
$data[201,2] = 0.08
$data[202,0] = "004612043"
$data[202,1] = "YURI"
$data[202,2] = 0.06
$data[203,0] = "004643880"
$data[203,1] = "GABRIEL"
$data[203,2] = 0.05
$data[204,0] = "004223226"
$data[204,1] = "YESHUA"
$data[204,2] = 0.04
$data[205,0] = "004281300"
$data[205,1] = "FRANKLIN"
$data[205,2] = 0.04
$data[206,0] = "005274028"
$data[206,1] = "RAFAEL"
$data[206,2] = 0.04
$data[207,0] = "004329229"
$data[207,1] = "GABRIEL"
$data[207,2] = 0.03
$data[208,0] = "004213373"
$data[208,1] = "ALEXANDRE"
$data[208,2] = 0.02
$data[209,0] = "004339183"
$data[209,1] = "JALISON"
$data[209,2] = 0.02
$data[210,0] = "000938440"
$data[210,1] = "BASE"
$data[210,2] = 0.01
$data[211,0] = "002878817"
$data[211,1] = "GUILHERME"
$data[211,2] = 0.01
$data[212,0] = "004400000"
$data[212,1] = "VILMA"
$data[212,2] = 0.01
$data[213,0] = ""
$data[213,1] = ""
$data[213,2] = ""
$data[214,0] = "Filtros aplicados:`nDataFim é (Em branco)`nnr_saldo_disponivel não é 0`nPosição é Posição D-1`nCARTEIRA não está em branco`ntela é Financeiro`nDataFim é (Em branco)`nNR_CONTA não está em branco`nTIPO_LANCAMENTO não é ED, ET ou Liquidação Doador"
$data[214,1] = ""
$data[214,2] = ""

# Conta (account number) values have leading zeros and must stay text,
# otherwise Excel would helpfully "fix" them into plain numbers.
$ws.Range("A2:A215").NumberFormat = "@"

$ws.Range("A1:C215").Value = $data
Write-Host "Export sheet rewritten: 215 rows"
